# Week 15 simulations update
# Applies appended per-game simulation numbers to the long space-separated
# strings on YDS / ST sheets, and refreshes the season-total rows on
# OFF / DEF / ST / TURNS / PEN sheets.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# YDS sheet - append new simulated play-by-play yardage values
# ---------------------------------------------------------------------
$ydsWs = $wb.Worksheets.Item("YDS")

$ydsWs.Range("B2").Value = $ydsWs.Range("B2").Value2 + " 2 1 3 3 2 5 1 6 15 1 2 13 -2 5 6 3 4 11 8 11 1 2 7 1 0 -2 1"
$ydsWs.Range("B3").Value = $ydsWs.Range("B3").Value2 + " -1 8 3 9 4 34 5 5 9 7 6 27 9 1 12 5 7 13 17 3 3 4"
$ydsWs.Range("C2").Value = $ydsWs.Range("C2").Value2 + " 5 1 4 5 4 3 -4 -1 7 2 18 4 8 5 0 11 13 21 5 3 2 1 1"
$ydsWs.Range("C3").Value = $ydsWs.Range("C3").Value2 + " 11 2 -1 5 3 5 4 9 4 20 10 32 4 8 -3 20 36 12 3 11 2 6 5 6 4 8 14"

# ---------------------------------------------------------------------
# OFF sheet - updated season totals
# ---------------------------------------------------------------------
$offWs = $wb.Worksheets.Item("OFF")

$offWs.Range("C2").Value = 175
$offWs.Range("D2").Value = 6
$offWs.Range("E2").Value = 9
$offWs.Range("F2").Value = 54
$offWs.Range("G2").Value = 48
$offWs.Range("H2").Value = 5
$offWs.Range("I2").Value = 4
$offWs.Range("J2").Value = 27
$offWs.Range("L2").Value = 208
$offWs.Range("M2").Value = 136
$offWs.Range("Q2").Value = 440

$offWs.Range("B3").Value = 12
$offWs.Range("C3").Value = 118
$offWs.Range("D3").Value = 4
$offWs.Range("F3").Value = 78
$offWs.Range("G3").Value = 21
$offWs.Range("H3").Value = 26
$offWs.Range("I3").Value = 36
$offWs.Range("J3").Value = 33

# ---------------------------------------------------------------------
# DEF sheet - updated season totals
# ---------------------------------------------------------------------
$defWs = $wb.Worksheets.Item("DEF")

$defWs.Range("B2").Value = 6
$defWs.Range("C2").Value = 146
$defWs.Range("D2").Value = 13
$defWs.Range("E2").Value = 10
$defWs.Range("F2").Value = 52
$defWs.Range("G2").Value = 33
$defWs.Range("I2").Value = 6
$defWs.Range("J2").Value = 22
$defWs.Range("L2").Value = 211
$defWs.Range("M2").Value = 135
$defWs.Range("O2").Value = 20
$defWs.Range("P2").Value = 12
$defWs.Range("Q2").Value = 403

$defWs.Range("C3").Value = 141
$defWs.Range("E3").Value = 15
$defWs.Range("F3").Value = 96
$defWs.Range("G3").Value = 28
$defWs.Range("H3").Value = 28
$defWs.Range("I3").Value = 55
$defWs.Range("J3").Value = 32
$defWs.Range("N3").Value = 20

# ---------------------------------------------------------------------
# ST sheet - updated season totals + appended simulation distributions
# ---------------------------------------------------------------------
$stWs = $wb.Worksheets.Item("ST")

$stWs.Range("B2").Value = 59
$stWs.Range("D2").Value = 45
$stWs.Range("F2").Value = 62
$stWs.Range("G2").Value = 60
$stWs.Range("L2").Value = 20
$stWs.Range("M2").Value = 10

$stWs.Range("B3").Value = 26

$stWs.Range("B4").Value = $stWs.Range("B4").Value2 + " 67 55 65"
$stWs.Range("B5").Value = $stWs.Range("B5").Value2 + " 19 12 23"
$stWs.Range("B6").Value = $stWs.Range("B6").Value2 + " 21 0 13 6"
$stWs.Range("D3").Value = $stWs.Range("D3").Value2 + " 38 35 39 41 43"
$stWs.Range("D4").Value = $stWs.Range("D4").Value2 + " 0 0 0 0 0"
$stWs.Range("D5").Value = $stWs.Range("D5").Value2 + " 8 1 6"

# ---------------------------------------------------------------------
# TURNS sheet - updated season totals
# ---------------------------------------------------------------------
$turnsWs = $wb.Worksheets.Item("TURNS")

$turnsWs.Range("B2").Value = 5
$turnsWs.Range("E2").Value = 7
$turnsWs.Range("D3").Value = 5

# ---------------------------------------------------------------------
# PEN sheet - updated season totals
# ---------------------------------------------------------------------
$penWs = $wb.Worksheets.Item("PEN")

$penWs.Range("B3").Value = 23
